# Mise à jour de l'application
# Adds a new wellness-survey session (2025-08-19, serial 45888) covering 17
# players, plus the new player name "Yoann Martelat" introduced by this
# session, as new rows 154-170 below the existing data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data as it appears in columns A-I:
#   row, playerName, Volume(C), Intensite(D), Fatigue(E), Douleur(F),
#   LocalisationDouleur(G, may be blank), Plaisir(H)
$rows = @(
    @(154, "Amir Etien",       55, 5, 6, 0, "",          6),
    @(155, "Yoann Martelat",   55, 5, 6, 5, "Genou",     7),
    @(156, "Omar Benyounes",   55, 5, 1, 0, "",          7),
    @(157, "Yanis Berrached",  55, 5, 7, 0, "",          0),
    @(158, "Romain Thunet",    55, 5, 5, 0, "",          2),
    @(159, "Maé Clavel",       55, 3, 2, 5, "Cheville",  6),
    @(160, "Yoan Zouma",       55, 3, 4, 4, "Adducteur",  7),
    @(161, "Levy Ndoutoume",   55, 3, 4, 3, "Ischio",    7),
    @(162, "Hedi Nasri",       55, 4, 3, 0, "",          8),
    @(163, "Ilan Ihaddadene",  55, 5, 3, 0, "",          5),
    @(164, "Emmanuel Valey",   55, 5, 5, 0, "",          5),
    @(165, "Karahali Souaré",  55, 5, 6, 6, "Cheville",  7),
    @(166, "Naim Dhib",        55, 4, 4, 0, "",          3),
    @(167, "Amir Kherrab",     55, 5, 6, 0, "",          8),
    @(168, "Wael Fareh",       55, 7, 5, 0, "",          8),
    @(169, "Mattheo Haon",     55, 8, 6, 0, "",          7),
    @(170, "Sofiane Belle",    55, 4, 3, 0, "",          5)
)

# Template rows already present in the sheet that carry the right cell
# styles: row 2 has an empty "Localisation douleur" (style used for blank
# G cells), row 153 has text in that column (style used when G is filled).
$blankTemplate = $ws.Range("A2:I2")
$textTemplate  = $ws.Range("A153:I153")

foreach ($r in $rows) {
    $rowNum = $r[0]
    $destRow = $ws.Range("A" + $rowNum + ":I" + $rowNum)

    if ($r[6] -eq "") {
        $blankTemplate.Copy()
    } else {
        $textTemplate.Copy()
    }
    $destRow.PasteSpecial(-4122)

    $ws.Range("A" + $rowNum).Value = 45888
    $ws.Range("B" + $rowNum).Value = $r[1]
    $ws.Range("C" + $rowNum).Value = $r[2]
    $ws.Range("D" + $rowNum).Value = $r[3]
    $ws.Range("E" + $rowNum).Value = $r[4]
    $ws.Range("F" + $rowNum).Value = $r[5]
    if ($r[6] -ne "") {
        $ws.Range("G" + $rowNum).Value = $r[6]
    }
    $ws.Range("H" + $rowNum).Value = $r[7]
}

# Assign the "Charge" formula to the whole new block in one shot so the
# engine stores it as a single shared formula (mirrors how row 131:153's
# I131*D131 pattern was already stored, just extended for the new rows).
$ws.Range("I154:I170").Formula = "=C154*D154"

$excel.CutCopyMode = $false

# Move the selection the way it ended up after the edit session.
$ws.Range("D174").Select()
